$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for the refreshed crypto symbol data.
# Cells are stored as text, so the number format is set to Text ("@") before
# assigning the value to prevent Excel from reinterpreting the string as a number.
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "303.23"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "5.23%"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "31.97"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "9.84%"
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "5.268"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "-0.19%"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.07498"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "6.93%"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "7.854"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "5.43%"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "3.816"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "7.29%"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "1.489"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "6.61%"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.9195"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "1.79%"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.1689"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "5.21%"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07903"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "4.82%"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.08019"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "3.80%"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.03041"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "4.31%"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.09895"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.001501"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "-4.52%"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.04603"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "1.75%"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.006560"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "8.23%"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "-0.68%"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "2.230"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "-0.07%"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "1.86%"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.1335"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "-0.73%"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.503"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "12.44%"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "1.43%"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.001215"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "0.58%"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.004447"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "6.95%"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.0001398"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "19.78%"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.0001933"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01719"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "2,534.35%"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.04480"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.006929"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "-0.39%"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "7.83%"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.002216"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "7.36%"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.01275"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "9.89%"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.00006148"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "4.84%"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.866"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "-3.30%"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.01497"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "15.34%"
